$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.568.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.578.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.44%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3709"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3353"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07496"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.959"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.581.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06781"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.423"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("E23").Value = "  -0.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.543.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.405"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.602"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.025"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.754.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.063"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.205"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.014"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.709"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08332"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02467"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2300"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.437"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06404"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.29%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6363"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.44%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6232"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.775"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.068"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.223"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07280"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "

Write-Host "Applied all cryptos updates"